$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITI")

# Mark every incident in the ITI backlog (rows 2-33, column I / "Status") as
# resolved.
$rng = $ws.Range("I2:I33")
$rng.Value = "Resolvido"

# A couple of rows (29 & 30) previously carried a different base style (the
# bold "Calibri" cell format) than the rest of the column. Re-apply the
# standard cell format used by the rest of the column (copied from I2) so the
# whole I2:I33 block shares one consistent style before the highlight fill is
# applied.
$ws.Range("I2").Copy()
$ws.Range("I29:I30").PasteSpecial(-4122)

# Highlight the whole status column with a yellow fill to call out the bulk
# update.
$rng.Interior.Color = 65535

# Leave the selection on the last couple of rows that were touched, matching
# where the user's cursor ended up after the bulk edit.
$ws.Activate()
$ws.Range("I32:I33").Select()
